# Generate Report for Handoff
# Inserts a new "12bb4fc5-9b48-497c-abef-af39a07043d4" handoff row above the
# existing "e13fe453-c82c-4c1e-97c3-ec00354eff20" row on every sheet.

$wb = $excel.ActiveWorkbook

$oldGuid = "e13fe453-c82c-4c1e-97c3-ec00354eff20"
$newGuid = "12bb4fc5-9b48-497c-abef-af39a07043d4"
$xlfHash = "d952076505a9ea68ca099dd13b9c3bc928acd92f"
$oldXlfHash = "990eaf09a03d5c07917fa1edfd5b2247ec1bb42e"

# ---------------------------------------------------------------------------
# Sheet "Overview": columns File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# The existing hyperlink is anchored to the row that is about to shift down,
# so drop it before the insert (it does not follow the shift automatically).
$wsOverview.Range("A2").Hyperlinks.Delete()

# Push the current data row down to row 3, leaving a blank row 2 for the new
# handoff entry.
$wsOverview.Rows.Item(2).Insert()

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("D2").Value = "2016-03-24 06:41:57"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/900e5016665140c33d96e82d1e4219c46c34f841/e2e/$newGuid.md", "", "", "$newGuid.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/900e5016665140c33d96e82d1e4219c46c34f841/e2e/$oldGuid.md", "", "", "$oldGuid.md")
$wsOverview.Range("A2").Style = "HyperLink"
$wsOverview.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheets "zh-cn" / "de-de": identical column layout, only locale-specific
# text differs.
# ---------------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Locale = "zh-cn"; NewHandoff = "2016-03-24 06:41:53"; OldHandoff = "2016-03-24 06:41:23"; NewHandbackUrlSeg = "53f43154bfdf11339f2fe75bedc6847a681f099c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn" },
    @{ Sheet = "de-de"; Locale = "de-de"; NewHandoff = "2016-03-24 06:41:57"; OldHandoff = "2016-03-24 06:41:27"; NewHandbackUrlSeg = "e52a96237dd251b9a54ee88eef490e713a691fd5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de" }
)

foreach ($info in $locales) {
    $ws = $wb.Worksheets.Item($info.Sheet)
    $locale = $info.Locale

    # Drop the hyperlinks that live on the row that is about to shift down.
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("D2").Hyperlinks.Delete()

    $ws.Rows.Item(2).Insert()

    $newMd = "$newGuid.md"
    $newXlf = "$newGuid.$xlfHash.$locale.xlf"
    $oldMd = "$oldGuid.md"
    $oldXlf = "$oldGuid.$oldXlfHash.$locale.xlf"

    $ws.Range("A2").Value = $newMd
    $ws.Range("B2").Value = ".md"
    $ws.Range("C2").Value = "Ready for handoff"
    $ws.Range("D2").Value = $newXlf
    $ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("E2").Value = $info.NewHandoff
    $ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("J2").Value = "Include"

    $ws.Range("A3").Value = $oldMd
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $oldXlf
    $ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("E3").Value = $info.OldHandoff
    $ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("J3").Value = "Include"

    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/900e5016665140c33d96e82d1e4219c46c34f841/e2e/$newMd", "", "", $newMd)
    $ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($info.NewHandbackUrlSeg)/ci/ht/$newXlf", "", "", $newXlf)
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/900e5016665140c33d96e82d1e4219c46c34f841/e2e/$oldMd", "", "", $oldMd)
    $ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($info.NewHandbackUrlSeg)/ci/ht/$oldXlf", "", "", $oldXlf)

    $ws.Range("A2").Style = "HyperLink"
    $ws.Range("D2").Style = "HyperLink"
    $ws.Range("A3").Style = "HyperLink"
    $ws.Range("D3").Style = "HyperLink"
}
